$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 170061253.5450011
$ws.Cells.Item(2, 3).Value = -402048568.8690004
$ws.Cells.Item(2, 4).Value = 53
$ws.Cells.Item(2, 5).Value = "T"
$ws.Cells.Item(3, 2).Value = 215751600.2529267
$ws.Cells.Item(3, 3).Value = -466156325.3815762
$ws.Cells.Item(3, 4).Value = 53
$ws.Cells.Item(3, 5).Value = "T"
$ws.Cells.Item(4, 2).Value = 190157095.1846463
$ws.Cells.Item(4, 3).Value = -415950649.4723316
$ws.Cells.Item(4, 4).Value = 53
$ws.Cells.Item(4, 5).Value = "T"
$ws.Cells.Item(5, 2).Value = 170860992.9712256
$ws.Cells.Item(5, 3).Value = -405155427.0714723
$ws.Cells.Item(5, 4).Value = 53
$ws.Cells.Item(5, 5).Value = "T"
$ws.Cells.Item(6, 2).Value = 170667185.431315
$ws.Cells.Item(6, 3).Value = -405104104.4883425
$ws.Cells.Item(6, 4).Value = 53
$ws.Cells.Item(6, 5).Value = "T"
$ws.Cells.Item(7, 2).Value = 198148318.8949276
$ws.Cells.Item(7, 3).Value = -446478395.0038871
$ws.Cells.Item(7, 4).Value = 53
$ws.Cells.Item(7, 5).Value = "T"
$ws.Cells.Item(8, 2).Value = 170922082.497805
$ws.Cells.Item(8, 3).Value = -404787896.7358894
$ws.Cells.Item(8, 4).Value = 53
$ws.Cells.Item(8, 5).Value = "T"
$ws.Cells.Item(9, 2).Value = 190811154.076052
$ws.Cells.Item(9, 3).Value = -418967515.8208964
$ws.Cells.Item(9, 4).Value = 53
$ws.Cells.Item(9, 5).Value = "T"
$ws.Cells.Item(10, 2).Value = 167063073.7771204
$ws.Cells.Item(10, 3).Value = -399736333.0616241
$ws.Cells.Item(10, 4).Value = 53
$ws.Cells.Item(10, 5).Value = "T"
$ws.Cells.Item(11, 2).Value = 191182675.188347
$ws.Cells.Item(11, 3).Value = -417313437.041097
$ws.Cells.Item(11, 4).Value = 53
$ws.Cells.Item(11, 5).Value = "T"
$ws.Cells.Item(12, 2).Value = 226125332.4069945
$ws.Cells.Item(12, 3).Value = -465476449.0371869
$ws.Cells.Item(12, 4).Value = 53
$ws.Cells.Item(12, 5).Value = "T"
$ws.Cells.Item(13, 2).Value = 207534924.9715279
$ws.Cells.Item(13, 3).Value = -448996303.2290542
$ws.Cells.Item(13, 4).Value = 53
$ws.Cells.Item(13, 5).Value = "T"
$ws.Cells.Item(14, 2).Value = 165417181.4159042
$ws.Cells.Item(14, 3).Value = -393879767.4921157
$ws.Cells.Item(14, 4).Value = 53
$ws.Cells.Item(14, 5).Value = "T"
$ws.Cells.Item(15, 2).Value = 179797278.9888872
$ws.Cells.Item(15, 3).Value = -385339572.3715404
$ws.Cells.Item(15, 4).Value = 53
$ws.Cells.Item(15, 5).Value = "T"
$ws.Cells.Item(16, 2).Value = 227210100.4166648
$ws.Cells.Item(16, 3).Value = -464887206.9620859
$ws.Cells.Item(16, 4).Value = 53
$ws.Cells.Item(16, 5).Value = "T"
$ws.Cells.Item(17, 2).Value = 212055122.9125671
$ws.Cells.Item(17, 3).Value = -453803865.4495928
$ws.Cells.Item(17, 4).Value = 53
$ws.Cells.Item(17, 5).Value = "T"
$ws.Cells.Item(18, 2).Value = 219066152.3479342
$ws.Cells.Item(18, 3).Value = -470844080.385905
$ws.Cells.Item(18, 4).Value = 53
$ws.Cells.Item(18, 5).Value = "T"
$ws.Cells.Item(19, 2).Value = 229520900.4912613
$ws.Cells.Item(19, 3).Value = -461115703.439372
$ws.Cells.Item(19, 4).Value = 53
$ws.Cells.Item(19, 5).Value = "T"
$ws.Cells.Item(20, 2).Value = 170189121.3721983
$ws.Cells.Item(20, 3).Value = -402456658.3516716
$ws.Cells.Item(20, 4).Value = 53
$ws.Cells.Item(20, 5).Value = "T"
$ws.Cells.Item(21, 2).Value = 232646418.6580028
$ws.Cells.Item(21, 3).Value = -467100958.0807744
$ws.Cells.Item(21, 4).Value = 53
$ws.Cells.Item(21, 5).Value = "T"
$ws.Cells.Item(22, 2).Value = 221685799.896555
$ws.Cells.Item(22, 3).Value = -450486429.1942855
$ws.Cells.Item(22, 4).Value = 53
$ws.Cells.Item(22, 5).Value = "T"
$ws.Cells.Item(23, 2).Value = 170044434.1339617
$ws.Cells.Item(23, 3).Value = -403681533.2859867
$ws.Cells.Item(23, 4).Value = 53
$ws.Cells.Item(23, 5).Value = "T"
$ws.Cells.Item(24, 2).Value = 170455328.0702242
$ws.Cells.Item(24, 3).Value = -404423592.289158
$ws.Cells.Item(24, 4).Value = 53
$ws.Cells.Item(24, 5).Value = "T"
$ws.Cells.Item(25, 2).Value = 220587369.5919776
$ws.Cells.Item(25, 3).Value = -459065039.3759356
$ws.Cells.Item(25, 4).Value = 53
$ws.Cells.Item(25, 5).Value = "T"
$ws.Cells.Item(26, 2).Value = 204294828.1179895
$ws.Cells.Item(26, 3).Value = -439678649.7171889
$ws.Cells.Item(26, 4).Value = 53
$ws.Cells.Item(26, 5).Value = "T"
$ws.Cells.Item(27, 2).Value = 169897527.0426158
$ws.Cells.Item(27, 3).Value = -402644194.6175289
$ws.Cells.Item(27, 4).Value = 53
$ws.Cells.Item(27, 5).Value = "T"
$ws.Cells.Item(28, 2).Value = 209752092.0299546
$ws.Cells.Item(28, 3).Value = -429279032.9395711
$ws.Cells.Item(28, 4).Value = 53
$ws.Cells.Item(28, 5).Value = "T"
$ws.Cells.Item(29, 2).Value = 181148880.3769667
$ws.Cells.Item(29, 3).Value = -393871673.9209508
$ws.Cells.Item(29, 4).Value = 53
$ws.Cells.Item(29, 5).Value = "T"
$ws.Cells.Item(30, 2).Value = 223592087.1595653
$ws.Cells.Item(30, 3).Value = -454483925.6360297
$ws.Cells.Item(30, 4).Value = 53
$ws.Cells.Item(30, 5).Value = "T"
$ws.Cells.Item(31, 2).Value = 220591388.8574566
$ws.Cells.Item(31, 3).Value = -459070734.1471248
$ws.Cells.Item(31, 4).Value = 53
$ws.Cells.Item(31, 5).Value = "T"
$ws.Cells.Item(32, 2).Value = 169837139.3715255
$ws.Cells.Item(32, 3).Value = -403010049.6333153
$ws.Cells.Item(32, 4).Value = 53
$ws.Cells.Item(32, 5).Value = "T"
$ws.Cells.Item(33, 2).Value = 217075885.120332
$ws.Cells.Item(33, 3).Value = -469975118.9094672
$ws.Cells.Item(33, 4).Value = 53
$ws.Cells.Item(33, 5).Value = "T"
$ws.Cells.Item(34, 2).Value = 170221339.6758617
$ws.Cells.Item(34, 3).Value = -404239109.246271
$ws.Cells.Item(34, 4).Value = 53
$ws.Cells.Item(34, 5).Value = "T"
$ws.Cells.Item(35, 2).Value = 228674745.6979045
$ws.Cells.Item(35, 3).Value = -459433627.7697359
$ws.Cells.Item(35, 4).Value = 53
$ws.Cells.Item(35, 5).Value = "T"
$ws.Cells.Item(36, 2).Value = 226511925.3198773
$ws.Cells.Item(36, 3).Value = -464173405.9727495
$ws.Cells.Item(36, 4).Value = 53
$ws.Cells.Item(36, 5).Value = "T"
$ws.Cells.Item(37, 2).Value = 182849971.5594774
$ws.Cells.Item(37, 3).Value = -397265861.4077961
$ws.Cells.Item(37, 4).Value = 53
$ws.Cells.Item(37, 5).Value = "T"
$ws.Cells.Item(38, 2).Value = 185605712.1449671
$ws.Cells.Item(38, 3).Value = -399507878.409721
$ws.Cells.Item(38, 4).Value = 53
$ws.Cells.Item(38, 5).Value = "T"
$ws.Cells.Item(39, 2).Value = 165158688.1655019
$ws.Cells.Item(39, 3).Value = -397520250.4713072
$ws.Cells.Item(39, 4).Value = 53
$ws.Cells.Item(39, 5).Value = "T"
$ws.Cells.Item(40, 2).Value = 169255388.5440889
$ws.Cells.Item(40, 3).Value = -400598228.2076799
$ws.Cells.Item(40, 4).Value = 53
$ws.Cells.Item(40, 5).Value = "T"
$ws.Cells.Item(41, 2).Value = 220469595.8216356
$ws.Cells.Item(41, 3).Value = -458743532.120712
$ws.Cells.Item(41, 4).Value = 53
$ws.Cells.Item(41, 5).Value = "T"
$ws.Cells.Item(42, 2).Value = 168088018.066688
$ws.Cells.Item(42, 3).Value = -400784453.2981902
$ws.Cells.Item(42, 4).Value = 53
$ws.Cells.Item(42, 5).Value = "T"
$ws.Cells.Item(43, 2).Value = 167778416.8765537
$ws.Cells.Item(43, 3).Value = -396994263.684976
$ws.Cells.Item(43, 4).Value = 53
$ws.Cells.Item(43, 5).Value = "T"
$ws.Cells.Item(44, 2).Value = 189893264.5149418
$ws.Cells.Item(44, 3).Value = -415728052.2039195
$ws.Cells.Item(44, 4).Value = 53
$ws.Cells.Item(44, 5).Value = "T"
$ws.Cells.Item(45, 2).Value = 161294915.6318805
$ws.Cells.Item(45, 3).Value = -386106025.2255547
$ws.Cells.Item(45, 4).Value = 53
$ws.Cells.Item(45, 5).Value = "T"
$ws.Cells.Item(46, 2).Value = 170989661.2193642
$ws.Cells.Item(46, 3).Value = -405566621.7426242
$ws.Cells.Item(46, 4).Value = 53
$ws.Cells.Item(46, 5).Value = "T"
$ws.Cells.Item(47, 2).Value = 167842617.8236445
$ws.Cells.Item(47, 3).Value = -401689953.4288404
$ws.Cells.Item(47, 4).Value = 53
$ws.Cells.Item(47, 5).Value = "T"
$ws.Cells.Item(48, 2).Value = 168348387.5857986
$ws.Cells.Item(48, 3).Value = -398256561.2255399
$ws.Cells.Item(48, 4).Value = 53
$ws.Cells.Item(48, 5).Value = "T"
$ws.Cells.Item(49, 2).Value = 225037801.297646
$ws.Cells.Item(49, 3).Value = -461912078.2331657
$ws.Cells.Item(49, 4).Value = 53
$ws.Cells.Item(49, 5).Value = "T"
$ws.Cells.Item(50, 2).Value = 170869051.6025432
$ws.Cells.Item(50, 3).Value = -403499299.370443
$ws.Cells.Item(50, 4).Value = 53
$ws.Cells.Item(50, 5).Value = "T"
$ws.Cells.Item(51, 2).Value = 193017470.4955127
$ws.Cells.Item(51, 3).Value = -401785101.312815
$ws.Cells.Item(51, 4).Value = 53
$ws.Cells.Item(51, 5).Value = "T"
$ws.Cells.Item(52, 2).Value = 199949949.5776266
$ws.Cells.Item(52, 3).Value = -448945059.420067
$ws.Cells.Item(52, 4).Value = 53
$ws.Cells.Item(52, 5).Value = "T"
$ws.Cells.Item(53, 2).Value = 212560183.8336233
$ws.Cells.Item(53, 3).Value = -460556023.6364397
$ws.Cells.Item(53, 4).Value = 53
$ws.Cells.Item(53, 5).Value = "T"
$ws.Cells.Item(54, 2).Value = 161676484.2880752
$ws.Cells.Item(54, 3).Value = -386240761.038307
$ws.Cells.Item(54, 4).Value = 53
$ws.Cells.Item(54, 5).Value = "T"
$ws.Cells.Item(55, 2).Value = 169788566.1146641
$ws.Cells.Item(55, 3).Value = -401736005.8829324
$ws.Cells.Item(55, 4).Value = 53
$ws.Cells.Item(55, 5).Value = "T"
$ws.Cells.Item(56, 2).Value = 211508905.0259657
$ws.Cells.Item(56, 3).Value = -452247932.9809692
$ws.Cells.Item(56, 4).Value = 53
$ws.Cells.Item(56, 5).Value = "T"
$ws.Cells.Item(57, 2).Value = 159944451.6321533
$ws.Cells.Item(57, 3).Value = -388908680.488093
$ws.Cells.Item(57, 4).Value = 53
$ws.Cells.Item(57, 5).Value = "T"
$ws.Cells.Item(58, 2).Value = 201934275.4551795
$ws.Cells.Item(58, 3).Value = -451332610.6898304
$ws.Cells.Item(58, 4).Value = 53
$ws.Cells.Item(58, 5).Value = "T"
$ws.Cells.Item(59, 2).Value = 184732978.0088556
$ws.Cells.Item(59, 3).Value = -397506742.2204202
$ws.Cells.Item(59, 4).Value = 53
$ws.Cells.Item(59, 5).Value = "T"
$ws.Cells.Item(60, 2).Value = 195705534.4388087
$ws.Cells.Item(60, 3).Value = -440888996.2268443
$ws.Cells.Item(60, 4).Value = 53
$ws.Cells.Item(60, 5).Value = "T"
$ws.Cells.Item(61, 2).Value = 220534136.2407041
$ws.Cells.Item(61, 3).Value = -482430540.6605094
$ws.Cells.Item(61, 4).Value = 53
$ws.Cells.Item(61, 5).Value = "T"
$ws.Cells.Item(62, 2).Value = 225035389.0068487
$ws.Cells.Item(62, 3).Value = -461908760.2337005
$ws.Cells.Item(62, 4).Value = 53
$ws.Cells.Item(62, 5).Value = "T"
$ws.Cells.Item(63, 2).Value = 183320051.9313003
$ws.Cells.Item(63, 3).Value = -390873338.4612513
$ws.Cells.Item(63, 4).Value = 53
$ws.Cells.Item(63, 5).Value = "T"
$ws.Cells.Item(64, 2).Value = 198119290.4729356
$ws.Cells.Item(64, 3).Value = -445812024.3581349
$ws.Cells.Item(64, 4).Value = 53
$ws.Cells.Item(64, 5).Value = "T"
$ws.Cells.Item(65, 2).Value = 183136998.266023
$ws.Cells.Item(65, 3).Value = -425597138.2776437
$ws.Cells.Item(65, 4).Value = 53
$ws.Cells.Item(65, 5).Value = "T"
$ws.Cells.Item(66, 2).Value = 194155942.949601
$ws.Cells.Item(66, 3).Value = -401840817.6504443
$ws.Cells.Item(66, 4).Value = 53
$ws.Cells.Item(66, 5).Value = "T"
$ws.Cells.Item(67, 2).Value = 199951766.5902286
$ws.Cells.Item(67, 3).Value = -448948024.5963792
$ws.Cells.Item(67, 4).Value = 53
$ws.Cells.Item(67, 5).Value = "T"
$ws.Cells.Item(68, 2).Value = 181191745.9251307
$ws.Cells.Item(68, 3).Value = -392441938.8805223
$ws.Cells.Item(68, 4).Value = 53
$ws.Cells.Item(68, 5).Value = "T"
$ws.Cells.Item(69, 2).Value = 181465006.6044621
$ws.Cells.Item(69, 3).Value = -381565732.0586818
$ws.Cells.Item(69, 4).Value = 53
$ws.Cells.Item(69, 5).Value = "T"
$ws.Cells.Item(70, 2).Value = 190126419.3274445
$ws.Cells.Item(70, 3).Value = -405131540.1544722
$ws.Cells.Item(70, 4).Value = 53
$ws.Cells.Item(70, 5).Value = "T"
$ws.Cells.Item(71, 2).Value = 213293560.7974027
$ws.Cells.Item(71, 3).Value = -444046630.9603429
$ws.Cells.Item(71, 4).Value = 53
$ws.Cells.Item(71, 5).Value = "T"
$ws.Cells.Item(72, 2).Value = 170191838.1157665
$ws.Cells.Item(72, 3).Value = -402461543.5417274
$ws.Cells.Item(72, 4).Value = 53
$ws.Cells.Item(72, 5).Value = "T"
$ws.Cells.Item(73, 2).Value = 220588452.1677152
$ws.Cells.Item(73, 3).Value = -459066573.2519769
$ws.Cells.Item(73, 4).Value = 53
$ws.Cells.Item(73, 5).Value = "T"
$ws.Cells.Item(74, 2).Value = 195997045.5293688
$ws.Cells.Item(74, 3).Value = -431996651.364231
$ws.Cells.Item(74, 4).Value = 53
$ws.Cells.Item(74, 5).Value = "T"
$ws.Cells.Item(75, 2).Value = 168450500.9432539
$ws.Cells.Item(75, 3).Value = -401384793.0106061
$ws.Cells.Item(75, 4).Value = 53
$ws.Cells.Item(75, 5).Value = "T"
$ws.Cells.Item(76, 2).Value = 230175350.2550632
$ws.Cells.Item(76, 3).Value = -469965199.3425449
$ws.Cells.Item(76, 4).Value = 53
$ws.Cells.Item(76, 5).Value = "T"
$ws.Cells.Item(77, 2).Value = 208808591.8167622
$ws.Cells.Item(77, 3).Value = -437210044.3358359
$ws.Cells.Item(77, 4).Value = 53
$ws.Cells.Item(77, 5).Value = "T"
$ws.Cells.Item(78, 2).Value = 175243179.3873374
$ws.Cells.Item(78, 3).Value = -381724269.4209624
$ws.Cells.Item(78, 4).Value = 53
$ws.Cells.Item(78, 5).Value = "T"
$ws.Cells.Item(79, 2).Value = 159981689.5436021
$ws.Cells.Item(79, 3).Value = -389034951.2546384
$ws.Cells.Item(79, 4).Value = 53
$ws.Cells.Item(79, 5).Value = "T"
$ws.Cells.Item(80, 2).Value = 201940941.1117887
$ws.Cells.Item(80, 3).Value = -451343395.5166952
$ws.Cells.Item(80, 4).Value = 53
$ws.Cells.Item(80, 5).Value = "T"
$ws.Cells.Item(81, 2).Value = 199038252.9463307
$ws.Cells.Item(81, 3).Value = -447386185.6953286
$ws.Cells.Item(81, 4).Value = 53
$ws.Cells.Item(81, 5).Value = "T"
$ws.Cells.Item(82, 2).Value = 200536386.0293367
$ws.Cells.Item(82, 3).Value = -450110034.1736351
$ws.Cells.Item(82, 4).Value = 53
$ws.Cells.Item(82, 5).Value = "T"
$ws.Cells.Item(83, 2).Value = 205759380.4885022
$ws.Cells.Item(83, 3).Value = -423798822.2369605
$ws.Cells.Item(83, 4).Value = 53
$ws.Cells.Item(83, 5).Value = "T"
$ws.Cells.Item(84, 2).Value = 189088595.2522518
$ws.Cells.Item(84, 3).Value = -395445249.0775359
$ws.Cells.Item(84, 4).Value = 53
$ws.Cells.Item(84, 5).Value = "T"
$ws.Cells.Item(85, 2).Value = 223641556.3890886
$ws.Cells.Item(85, 3).Value = -459278078.5597609
$ws.Cells.Item(85, 4).Value = 53
$ws.Cells.Item(85, 5).Value = "T"
$ws.Cells.Item(86, 2).Value = 201021809.2568428
$ws.Cells.Item(86, 3).Value = -450970544.5014009
$ws.Cells.Item(86, 4).Value = 53
$ws.Cells.Item(86, 5).Value = "T"
$ws.Cells.Item(87, 2).Value = 180714862.6879153
$ws.Cells.Item(87, 3).Value = -386459896.4466831
$ws.Cells.Item(87, 4).Value = 53
$ws.Cells.Item(87, 5).Value = "T"
$ws.Cells.Item(88, 2).Value = 192332759.4216259
$ws.Cells.Item(88, 3).Value = -434773788.9701531
$ws.Cells.Item(88, 4).Value = 53
$ws.Cells.Item(88, 5).Value = "T"
$ws.Cells.Item(89, 2).Value = 160418897.5706784
$ws.Cells.Item(89, 3).Value = -389911002.6655887
$ws.Cells.Item(89, 4).Value = 53
$ws.Cells.Item(89, 5).Value = "T"
$ws.Cells.Item(90, 2).Value = 183326049.5143205
$ws.Cells.Item(90, 3).Value = -425033345.1795313
$ws.Cells.Item(90, 4).Value = 53
$ws.Cells.Item(90, 5).Value = "T"
$ws.Cells.Item(91, 2).Value = 159985281.4310923
$ws.Cells.Item(91, 3).Value = -389041687.233989
$ws.Cells.Item(91, 4).Value = 53
$ws.Cells.Item(91, 5).Value = "T"
$ws.Cells.Item(92, 2).Value = 178908355.1332637
$ws.Cells.Item(92, 3).Value = -404735969.3286093
$ws.Cells.Item(92, 4).Value = 53
$ws.Cells.Item(92, 5).Value = "T"
$ws.Cells.Item(93, 2).Value = 158953626.7530148
$ws.Cells.Item(93, 3).Value = -386767804.1345814
$ws.Cells.Item(93, 4).Value = 53
$ws.Cells.Item(93, 5).Value = "T"
$ws.Cells.Item(94, 2).Value = 159983783.1963968
$ws.Cells.Item(94, 3).Value = -389038877.5547693
$ws.Cells.Item(94, 4).Value = 53
$ws.Cells.Item(94, 5).Value = "T"
$ws.Cells.Item(95, 2).Value = 159031701.4770814
$ws.Cells.Item(95, 3).Value = -387026708.7799101
$ws.Cells.Item(95, 4).Value = 53
$ws.Cells.Item(95, 5).Value = "T"
$ws.Cells.Item(96, 2).Value = 191708627.3984575
$ws.Cells.Item(96, 3).Value = -410724560.5956822
$ws.Cells.Item(96, 4).Value = 53
$ws.Cells.Item(96, 5).Value = "T"
$ws.Cells.Item(97, 2).Value = 164829082.7846155
$ws.Cells.Item(97, 3).Value = -396448779.1607118
$ws.Cells.Item(97, 4).Value = 53
$ws.Cells.Item(97, 5).Value = "T"
$ws.Cells.Item(98, 2).Value = 191830326.5457348
$ws.Cells.Item(98, 3).Value = -420857416.8751051
$ws.Cells.Item(98, 4).Value = 53
$ws.Cells.Item(98, 5).Value = "T"
$ws.Cells.Item(99, 2).Value = 200745874.1465958
$ws.Cells.Item(99, 3).Value = -450730133.2302119
$ws.Cells.Item(99, 4).Value = 53
$ws.Cells.Item(99, 5).Value = "T"
$ws.Cells.Item(100, 2).Value = 195100893.162625
$ws.Cells.Item(100, 3).Value = -440230697.0622374
$ws.Cells.Item(100, 4).Value = 53
$ws.Cells.Item(100, 5).Value = "T"
$ws.Cells.Item(101, 2).Value = 225334535.0004603
$ws.Cells.Item(101, 3).Value = -455721416.1328683
$ws.Cells.Item(101, 4).Value = 53
$ws.Cells.Item(101, 5).Value = "T"
$ws.Cells.Item(102, 2).Value = 170066567.4634582
$ws.Cells.Item(102, 3).Value = -402058122.279879
$ws.Cells.Item(102, 4).Value = 53
$ws.Cells.Item(102, 5).Value = "T"
$ws.Cells.Item(103, 2).Value = 184677066.2989198
$ws.Cells.Item(103, 3).Value = -426343195.7846923
$ws.Cells.Item(103, 4).Value = 53
$ws.Cells.Item(103, 5).Value = "T"
$ws.Cells.Item(104, 2).Value = 190705854.829486
$ws.Cells.Item(104, 3).Value = -429290698.053645
$ws.Cells.Item(104, 4).Value = 53
$ws.Cells.Item(104, 5).Value = "T"
$ws.Cells.Item(105, 2).Value = 173950496.5190559
$ws.Cells.Item(105, 3).Value = -404282480.4928117
$ws.Cells.Item(105, 4).Value = 53
$ws.Cells.Item(105, 5).Value = "T"
$ws.Cells.Item(106, 2).Value = 236488082.5247939
$ws.Cells.Item(106, 3).Value = -502234364.1522337
$ws.Cells.Item(106, 4).Value = 53
$ws.Cells.Item(106, 5).Value = "T"
$ws.Cells.Item(107, 2).Value = 212723687.7217174
$ws.Cells.Item(107, 3).Value = -449317799.2371773
$ws.Cells.Item(107, 4).Value = 53
$ws.Cells.Item(107, 5).Value = "T"
$ws.Cells.Item(108, 2).Value = 169111069.225488
$ws.Cells.Item(108, 3).Value = -402939472.3102995
$ws.Cells.Item(108, 4).Value = 53
$ws.Cells.Item(108, 5).Value = "T"
$ws.Cells.Item(109, 2).Value = 198336644.9666787
$ws.Cells.Item(109, 3).Value = -411866837.4992051
$ws.Cells.Item(109, 4).Value = 53
$ws.Cells.Item(109, 5).Value = "T"
$ws.Cells.Item(110, 2).Value = 158546247.0333468
$ws.Cells.Item(110, 3).Value = -380475198.7775751
$ws.Cells.Item(110, 4).Value = 53
$ws.Cells.Item(110, 5).Value = "T"
$ws.Cells.Item(111, 2).Value = 189729394.5335976
$ws.Cells.Item(111, 3).Value = -403999810.6284639
$ws.Cells.Item(111, 4).Value = 53
$ws.Cells.Item(111, 5).Value = "T"
$ws.Cells.Item(112, 2).Value = 202072665.5388272
$ws.Cells.Item(112, 3).Value = -453520694.1747884
$ws.Cells.Item(112, 4).Value = 53
$ws.Cells.Item(112, 5).Value = "T"
$ws.Cells.Item(113, 2).Value = 161325373.0124752
$ws.Cells.Item(113, 3).Value = -385652928.7783066
$ws.Cells.Item(113, 4).Value = 53
$ws.Cells.Item(113, 5).Value = "T"
$ws.Cells.Item(114, 2).Value = 161054255.9680095
$ws.Cells.Item(114, 3).Value = -385324290.4075406
$ws.Cells.Item(114, 4).Value = 53
$ws.Cells.Item(114, 5).Value = "T"
$ws.Cells.Item(115, 2).Value = 228912456.7268519
$ws.Cells.Item(115, 3).Value = -492802897.0585838
$ws.Cells.Item(115, 4).Value = 53
$ws.Cells.Item(115, 5).Value = "T"
$ws.Cells.Item(116, 2).Value = 189903671.6324413
$ws.Cells.Item(116, 3).Value = -408201000.6104372
$ws.Cells.Item(116, 4).Value = 53
$ws.Cells.Item(116, 5).Value = "T"
$ws.Cells.Item(117, 2).Value = 190128098.1234068
$ws.Cells.Item(117, 3).Value = -415313302.3688121
$ws.Cells.Item(117, 4).Value = 53
$ws.Cells.Item(117, 5).Value = "T"
$ws.Cells.Item(118, 2).Value = 219741114.3635682
$ws.Cells.Item(118, 3).Value = -454396973.9048197
$ws.Cells.Item(118, 4).Value = 53
$ws.Cells.Item(118, 5).Value = "T"
$ws.Cells.Item(119, 2).Value = 191441505.7863299
$ws.Cells.Item(119, 3).Value = -427522610.3461448
$ws.Cells.Item(119, 4).Value = 53
$ws.Cells.Item(119, 5).Value = "T"
$ws.Cells.Item(120, 2).Value = 166744254.6523182
$ws.Cells.Item(120, 3).Value = -395359670.511667
$ws.Cells.Item(120, 4).Value = 53
$ws.Cells.Item(120, 5).Value = "T"
$ws.Cells.Item(121, 2).Value = 159982778.8322647
$ws.Cells.Item(121, 3).Value = -389036994.0382276
$ws.Cells.Item(121, 4).Value = 53
$ws.Cells.Item(121, 5).Value = "T"
$ws.Cells.Item(122, 2).Value = 159450693.4018503
$ws.Cells.Item(122, 3).Value = -382313786.5138384
$ws.Cells.Item(122, 4).Value = 53
$ws.Cells.Item(122, 5).Value = "T"
$ws.Cells.Item(123, 2).Value = 160581153.2411571
$ws.Cells.Item(123, 3).Value = -384341226.8542533
$ws.Cells.Item(123, 4).Value = 53
$ws.Cells.Item(123, 5).Value = "T"
$ws.Cells.Item(124, 2).Value = 204412644.2140465
$ws.Cells.Item(124, 3).Value = -432688320.1416128
$ws.Cells.Item(124, 4).Value = 53
$ws.Cells.Item(124, 5).Value = "T"
$ws.Cells.Item(125, 2).Value = 163617850.8680348
$ws.Cells.Item(125, 3).Value = -388648087.9597438
$ws.Cells.Item(125, 4).Value = 53
$ws.Cells.Item(125, 5).Value = "T"
$ws.Cells.Item(126, 2).Value = 188672409.5696528
$ws.Cells.Item(126, 3).Value = -412151951.1620351
$ws.Cells.Item(126, 4).Value = 53
$ws.Cells.Item(126, 5).Value = "T"
$ws.Cells.Item(127, 2).Value = 180662603.712813
$ws.Cells.Item(127, 3).Value = -380285375.9184469
$ws.Cells.Item(127, 4).Value = 53
$ws.Cells.Item(127, 5).Value = "T"
$ws.Cells.Item(128, 2).Value = 159981312.4340496
$ws.Cells.Item(128, 3).Value = -389034244.04576
$ws.Cells.Item(128, 4).Value = 53
$ws.Cells.Item(128, 5).Value = "T"
$ws.Cells.Item(129, 2).Value = 190553921.2513739
$ws.Cells.Item(129, 3).Value = -418754253.6420302
$ws.Cells.Item(129, 4).Value = 53
$ws.Cells.Item(129, 5).Value = "T"
$ws.Cells.Item(130, 2).Value = 164680010.3926166
$ws.Cells.Item(130, 3).Value = -379685424.6714844
$ws.Cells.Item(130, 4).Value = 53
$ws.Cells.Item(130, 5).Value = "T"
$ws.Cells.Item(131, 2).Value = 180693895.2109168
$ws.Cells.Item(131, 3).Value = -407424798.7042264
$ws.Cells.Item(131, 4).Value = 53
$ws.Cells.Item(131, 5).Value = "T"
$ws.Cells.Item(132, 2).Value = 162283396.6694638
$ws.Cells.Item(132, 3).Value = -386557522.8375871
$ws.Cells.Item(132, 4).Value = 53
$ws.Cells.Item(132, 5).Value = "T"
$ws.Cells.Item(133, 2).Value = 159593752.1135564
$ws.Cells.Item(133, 3).Value = -388307085.0363675
$ws.Cells.Item(133, 4).Value = 53
$ws.Cells.Item(133, 5).Value = "T"
$ws.Cells.Item(134, 2).Value = 183568970.4222474
$ws.Cells.Item(134, 3).Value = -418382782.426956
$ws.Cells.Item(134, 4).Value = 53
$ws.Cells.Item(134, 5).Value = "T"
$ws.Cells.Item(135, 2).Value = 193274830.4161268
$ws.Cells.Item(135, 3).Value = -400423351.3916775
$ws.Cells.Item(135, 4).Value = 53
$ws.Cells.Item(135, 5).Value = "T"
$ws.Cells.Item(136, 2).Value = 169039221.4966475
$ws.Cells.Item(136, 3).Value = -398792434.4795803
$ws.Cells.Item(136, 4).Value = 53
$ws.Cells.Item(136, 5).Value = "T"
$ws.Cells.Item(137, 2).Value = 187767923.9888766
$ws.Cells.Item(137, 3).Value = -412204218.9547121
$ws.Cells.Item(137, 4).Value = 53
$ws.Cells.Item(137, 5).Value = "T"
$ws.Cells.Item(138, 2).Value = 176614885.9983608
$ws.Cells.Item(138, 3).Value = -404281953.2558571
$ws.Cells.Item(138, 4).Value = 53
$ws.Cells.Item(138, 5).Value = "T"
$ws.Cells.Item(139, 2).Value = 214146793.1152762
$ws.Cells.Item(139, 3).Value = -473637174.7999482
$ws.Cells.Item(139, 4).Value = 53
$ws.Cells.Item(139, 5).Value = "T"
$ws.Cells.Item(140, 2).Value = 195694588.3217317
$ws.Cells.Item(140, 3).Value = -437953639.0172651
$ws.Cells.Item(140, 4).Value = 53
$ws.Cells.Item(140, 5).Value = "T"
$ws.Cells.Item(141, 2).Value = 161681400.0862358
$ws.Cells.Item(141, 3).Value = -386249778.728037
$ws.Cells.Item(141, 4).Value = 53
$ws.Cells.Item(141, 5).Value = "T"
$ws.Cells.Item(142, 2).Value = 199808924.882672
$ws.Cells.Item(142, 3).Value = -421352675.5362641
$ws.Cells.Item(142, 4).Value = 53
$ws.Cells.Item(142, 5).Value = "T"
$ws.Cells.Item(143, 2).Value = 187914817.4255285
$ws.Cells.Item(143, 3).Value = -411011705.6538569
$ws.Cells.Item(143, 4).Value = 53
$ws.Cells.Item(143, 5).Value = "T"
$ws.Cells.Item(144, 2).Value = 170030570.7160479
$ws.Cells.Item(144, 3).Value = -373916385.7070824
$ws.Cells.Item(144, 4).Value = 53
$ws.Cells.Item(144, 5).Value = "T"
$ws.Cells.Item(145, 2).Value = 218929775.7707652
$ws.Cells.Item(145, 3).Value = -440791059.8359208
$ws.Cells.Item(145, 4).Value = 53
$ws.Cells.Item(145, 5).Value = "T"
$ws.Cells.Item(146, 2).Value = 162708234.4723862
$ws.Cells.Item(146, 3).Value = -388472907.154898
$ws.Cells.Item(146, 4).Value = 53
$ws.Cells.Item(146, 5).Value = "T"
$ws.Cells.Item(147, 2).Value = 193098927.3384862
$ws.Cells.Item(147, 3).Value = -436513085.7480596
$ws.Cells.Item(147, 4).Value = 53
$ws.Cells.Item(147, 5).Value = "T"
$ws.Cells.Item(148, 2).Value = 230555666.8963441
$ws.Cells.Item(148, 3).Value = -459769917.9560463
$ws.Cells.Item(148, 4).Value = 53
$ws.Cells.Item(148, 5).Value = "T"
$ws.Cells.Item(149, 2).Value = 193727869.7420573
$ws.Cells.Item(149, 3).Value = -433795557.1452767
$ws.Cells.Item(149, 4).Value = 53
$ws.Cells.Item(149, 5).Value = "T"
$ws.Cells.Item(150, 2).Value = 172808170.7906583
$ws.Cells.Item(150, 3).Value = -405162410.6592178
$ws.Cells.Item(150, 4).Value = 53
$ws.Cells.Item(150, 5).Value = "T"
$ws.Cells.Item(151, 2).Value = 190025554.5393861
$ws.Cells.Item(151, 3).Value = -418302784.1040679
$ws.Cells.Item(151, 4).Value = 53
$ws.Cells.Item(151, 5).Value = "T"
$ws.Cells.Item(152, 2).Value = 233758818.4293233
$ws.Cells.Item(152, 3).Value = -497716248.5745634
$ws.Cells.Item(152, 4).Value = 53
$ws.Cells.Item(152, 5).Value = "T"
$ws.Cells.Item(153, 2).Value = 197999048.2211913
$ws.Cells.Item(153, 3).Value = -442524433.6819119
$ws.Cells.Item(153, 4).Value = 53
$ws.Cells.Item(153, 5).Value = "T"
$ws.Cells.Item(154, 2).Value = 170274010.9602531
$ws.Cells.Item(154, 3).Value = -401599067.683601
$ws.Cells.Item(154, 4).Value = 53
$ws.Cells.Item(154, 5).Value = "T"
$ws.Cells.Item(155, 2).Value = 217091423.600224
$ws.Cells.Item(155, 3).Value = -452864053.2432491
$ws.Cells.Item(155, 4).Value = 53
$ws.Cells.Item(155, 5).Value = "T"
$ws.Cells.Item(156, 2).Value = 166872681.8275874
$ws.Cells.Item(156, 3).Value = -394659380.2129357
$ws.Cells.Item(156, 4).Value = 53
$ws.Cells.Item(156, 5).Value = "T"
$ws.Cells.Item(157, 2).Value = 208841063.2337069
$ws.Cells.Item(157, 3).Value = -433398500.873895
$ws.Cells.Item(157, 4).Value = 53
$ws.Cells.Item(157, 5).Value = "T"
$ws.Cells.Item(158, 2).Value = 191697416.7709543
$ws.Cells.Item(158, 3).Value = -400171482.9756765
$ws.Cells.Item(158, 4).Value = 53
$ws.Cells.Item(158, 5).Value = "T"
$ws.Cells.Item(159, 2).Value = 163624090.056206
$ws.Cells.Item(159, 3).Value = -390320148.7865415
$ws.Cells.Item(159, 4).Value = 53
$ws.Cells.Item(159, 5).Value = "T"
$ws.Cells.Item(160, 2).Value = 230432261.7004661
$ws.Cells.Item(160, 3).Value = -459446216.7997958
$ws.Cells.Item(160, 4).Value = 53
$ws.Cells.Item(160, 5).Value = "T"
$ws.Cells.Item(161, 2).Value = 213054465.7445063
$ws.Cells.Item(161, 3).Value = -433863851.7112849
$ws.Cells.Item(161, 4).Value = 53
$ws.Cells.Item(161, 5).Value = "T"
$ws.Cells.Item(162, 2).Value = 192706308.6065789
$ws.Cells.Item(162, 3).Value = -435326852.917184
$ws.Cells.Item(162, 4).Value = 53
$ws.Cells.Item(162, 5).Value = "T"
$ws.Cells.Item(163, 2).Value = 208919250.3967201
$ws.Cells.Item(163, 3).Value = -435841706.9338781
$ws.Cells.Item(163, 4).Value = 53
$ws.Cells.Item(163, 5).Value = "T"
$ws.Cells.Item(164, 2).Value = 169634345.0394309
$ws.Cells.Item(164, 3).Value = -400688439.0056636
$ws.Cells.Item(164, 4).Value = 53
$ws.Cells.Item(164, 5).Value = "T"
$ws.Cells.Item(165, 2).Value = 214076451.9032896
$ws.Cells.Item(165, 3).Value = -434970040.9152226
$ws.Cells.Item(165, 4).Value = 53
$ws.Cells.Item(165, 5).Value = "T"
$ws.Cells.Item(166, 2).Value = 226014241.3440962
$ws.Cells.Item(166, 3).Value = -460462448.9135036
$ws.Cells.Item(166, 4).Value = 53
$ws.Cells.Item(166, 5).Value = "T"
$ws.Cells.Item(167, 2).Value = 169510653.6271445
$ws.Cells.Item(167, 3).Value = -400288688.0457355
$ws.Cells.Item(167, 4).Value = 53
$ws.Cells.Item(167, 5).Value = "T"
$ws.Cells.Item(168, 2).Value = 197544351.2653601
$ws.Cells.Item(168, 3).Value = -428644976.7766817
$ws.Cells.Item(168, 4).Value = 53
$ws.Cells.Item(168, 5).Value = "T"
$ws.Cells.Item(169, 2).Value = 225147600.0272056
$ws.Cells.Item(169, 3).Value = -446621699.3602381
$ws.Cells.Item(169, 4).Value = 53
$ws.Cells.Item(169, 5).Value = "T"
$ws.Cells.Item(170, 2).Value = 204452339.7224762
$ws.Cells.Item(170, 3).Value = -427274707.5131188
$ws.Cells.Item(170, 4).Value = 53
$ws.Cells.Item(170, 5).Value = "T"
$ws.Cells.Item(171, 2).Value = 211436154.0418527
$ws.Cells.Item(171, 3).Value = -432216727.6757636
$ws.Cells.Item(171, 4).Value = 53
$ws.Cells.Item(171, 5).Value = "T"
$ws.Cells.Item(172, 2).Value = 224428130.6488923
$ws.Cells.Item(172, 3).Value = -460835459.0009965
$ws.Cells.Item(172, 4).Value = 53
$ws.Cells.Item(172, 5).Value = "T"
$ws.Cells.Item(173, 2).Value = 230867323.3203529
$ws.Cells.Item(173, 3).Value = -459433945.2831934
$ws.Cells.Item(173, 4).Value = 53
$ws.Cells.Item(173, 5).Value = "T"
$ws.Cells.Item(174, 2).Value = 188099519.5105205
$ws.Cells.Item(174, 3).Value = -414814430.0133743
$ws.Cells.Item(174, 4).Value = 53
$ws.Cells.Item(174, 5).Value = "T"
$ws.Cells.Item(175, 2).Value = 226042122.9440573
$ws.Cells.Item(175, 3).Value = -448979737.9827173
$ws.Cells.Item(175, 4).Value = 53
$ws.Cells.Item(175, 5).Value = "T"
$ws.Cells.Item(176, 2).Value = 227661733.7457307
$ws.Cells.Item(176, 3).Value = -466707593.4775785
$ws.Cells.Item(176, 4).Value = 53
$ws.Cells.Item(176, 5).Value = "T"
$ws.Cells.Item(177, 2).Value = 189877414.6763963
$ws.Cells.Item(177, 3).Value = -414587537.4387551
$ws.Cells.Item(177, 4).Value = 53
$ws.Cells.Item(177, 5).Value = "T"
$ws.Cells.Item(178, 2).Value = 194700531.8312197
$ws.Cells.Item(178, 3).Value = -439028991.2625205
$ws.Cells.Item(178, 4).Value = 53
$ws.Cells.Item(178, 5).Value = "T"
$ws.Cells.Item(179, 2).Value = 193084589.4705866
$ws.Cells.Item(179, 3).Value = -435886815.5279012
$ws.Cells.Item(179, 4).Value = 53
$ws.Cells.Item(179, 5).Value = "T"
$ws.Cells.Item(180, 2).Value = 215805372.3812869
$ws.Cells.Item(180, 3).Value = -465711154.7927014
$ws.Cells.Item(180, 4).Value = 53
$ws.Cells.Item(180, 5).Value = "T"
